$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("VENTA MENSUAL")

# Widen column B (CLIENTE) from 36 to 56 character-units.
# Excel's ColumnWidth property is expressed in "characters" of the Normal
# style font, which is offset from the width stored in the worksheet XML
# by ~5/6 of a character (5 padding pixels / 6px max digit width for the
# default Calibri 11 font). Subtracting that offset here makes the value
# that ends up persisted in the file equal to exactly 56.
$ws.Columns.Item(2).ColumnWidth = 56 - 5/6

# A new client, "CONSTRUCCION, INGENIERIA Y TECNOLOGIA CONSTRUINTEC SAS",
# was added alphabetically right after "CARAVEDO PAZMIÑO  JAHAIRA PAMELA"
# (row 6) and before "DANIELA ELIZABETH BECERRA BECERRA" (previously row
# 7). Insert a new row at 7; this shifts all following rows (and the
# totals row) down by one, preserving their values and formatting.
$ws.Rows.Item(7).Insert()

# Populate the newly inserted row with the new client's data (all
# monthly/budget figures are 0, so the totals row is unaffected).
$ws.Range("A7").Value = "OFICINA-CATAECSA"
$ws.Range("B7").Value = "CONSTRUCCION, INGENIERIA Y TECNOLOGIA CONSTRUINTEC SAS"
$ws.Range("C7").Value = 0
$ws.Range("D7").Value = 0
$ws.Range("E7").Value = 0
$ws.Range("F7").Value = 0
$ws.Range("G7").Value = 0
